# This script applies the cell-value updates to Sheet1 as described in the
# commit diff. Every change below targets a single numeric cell in the
# player statistics table (rows 2-101, columns B-N). Values are written
# using $ws.Range(...).Value so the resulting cached <v> in the OOXML
# matches the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("C32").Value = 0
$ws.Range("H43").Value = 0.058
$ws.Range("C46").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H62").Value = 0.06
$ws.Range("I62").Value = 0.05899999999999994
$ws.Range("J65").Value = 0
$ws.Range("G68").Value = 0.03000000000000003
$ws.Range("I68").Value = 0.02800000000000002
$ws.Range("D69").Value = 0.07500000000000001
$ws.Range("B71").Value = 0.003000000000000003
$ws.Range("E71").Value = -0.01300000000000001
$ws.Range("F77").Value = -0.006000000000000005
$ws.Range("H77").Value = 0.03299999999999997
$ws.Range("J79").Value = 0.02300000000000002
$ws.Range("F80").Value = 0.03700000000000003
$ws.Range("G81").Value = 0.09799999999999998
$ws.Range("B82").Value = 0
$ws.Range("C82").Value = 0
$ws.Range("J84").Value = 0.139
$ws.Range("E87").Value = -0.0129999999999999
$ws.Range("H87").Value = 0.02000000000000002
$ws.Range("I88").Value = 0.136
$ws.Range("G89").Value = 0.03700000000000003
$ws.Range("J90").Value = 0.04699999999999999
$ws.Range("B91").Value = -0.03700000000000003
$ws.Range("H91").Value = 0.04599999999999999
$ws.Range("H93").Value = 0.07200000000000001
$ws.Range("C95").Value = -0.0169999999999999
$ws.Range("L95").Value = -0.03899999999999998
$ws.Range("G97").Value = 0.1080000000000001
$ws.Range("L97").Value = 0.02200000000000002
$ws.Range("B98").Value = -0.02200000000000002
$ws.Range("K98").Value = -0.113
$ws.Range("C99").Value = 0.03100000000000003
$ws.Range("H99").Value = 0.04899999999999999
$ws.Range("C100").Value = -0.05000000000000004
$ws.Range("E100").Value = -0.03000000000000003

$wb.Save()
